# Fixed #476 Moving from Apache POI 4.1.0 to 5.2.3.
#
# The POI upgrade changed how run properties are (re)serialised for the
# generated table: ST_OnOff booleans are now written as "on"/"off"
# instead of "true"/"false", and <w:sz> is emitted first in <w:rPr>
# (ahead of <w:rFonts>/<w:b>/<w:i>/<w:strike>/<w:color>). Re-assert the
# same run formatting (Candara 11pt, black, not bold/italic/struck —
# except the header "Total" cell, which stays bold) on every data cell
# of the generated table so the properties are written back out with
# that layout.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function New-CellFragment([string]$text, [bool]$bold) {
    $bVal = "off"
    if ($bold) { $bVal = "on" }
    $escText = $text -replace '&','&amp;' -replace '<','&lt;' -replace '>','&gt;'
    $xml = @"
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:before="0" w:after="0"/></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:rFonts w:ascii="Candara" w:hAnsi="Candara" w:cs="Candara" w:eastAsia="Candara"/><w:b w:val="$bVal"/><w:i w:val="off"/><w:strike w:val="off"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">$escText</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
    return $xml
}

# row -> list of (text, bold) per column, matching the generated table.
$rows = @(
    @(@("Item", $false), @("Price", $false), @("Quantity", $false), @("Total", $true)),
    @(@("Apple", $false), @("1,76 €", $false), @("23", $false), @("40,48 €", $false)),
    @(@("Orange", $false), @("2,12 €", $false), @("13", $false), @("27,56 €", $false)),
    @(@("Banana", $false), @("1,99 €", $false), @("45", $false), @("89,55 €", $false)),
    @($null, $null, @("Total", $false), @("157,59 €", $false))
)

for ($rowIdx = 0; $rowIdx -lt $rows.Count; $rowIdx++) {
    $row = $rows[$rowIdx]
    for ($colIdx = 0; $colIdx -lt $row.Count; $colIdx++) {
        $entry = $row[$colIdx]
        if ($null -eq $entry) { continue }

        $cell = $t.Cell($rowIdx + 1, $colIdx + 1)
        $cellRange = $cell.Range
        # Exclude the trailing end-of-cell marker so InsertXML only
        # touches the paragraph that already holds the run.
        $target = $d.Range($cellRange.Start, $cellRange.End - 1)

        $frag = New-CellFragment $entry[0] $entry[1]
        [void]$target.InsertXML($frag)
    }
}


